$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2600.5715
$ws.Range("I86").Value = 1296.6086
$ws.Range("J86").Value = 8598.799999999999
$ws.Range("K86").Value = 1296.6086
$ws.Range("L86").Value = 8598.799999999999
$ws.Range("M86").Value = -173.6086
$ws.Range("N86").Value = -10844.8
$ws.Range("H89").Value = 2600.5715
$ws.Range("I89").Value = 1296.6086
$ws.Range("J89").Value = 8598.799999999999
$ws.Range("K89").Value = 6483.043
$ws.Range("L89").Value = 42994
$ws.Range("M89").Value = -867.0429999999997
$ws.Range("N89").Value = -54226
$ws.Range("H113").Value = 57398.332
$ws.Range("I113").Value = 334668.34
$ws.Range("J113").Value = 1944.3334
$ws.Range("K113").Value = 334668.34
$ws.Range("L113").Value = 1944.3334
$ws.Range("M113").Value = -331414.34
$ws.Range("N113").Value = -8452.3334
$ws.Range("H116").Value = 2101.25
$ws.Range("I116").Value = 2101.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2101.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1340.75
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 4721266
$ws.Range("I132").Value = 4812047
$ws.Range("J132").Value = 650
$ws.Range("K132").Value = 14436141
$ws.Range("L132").Value = 1950
$ws.Range("M132").Value = -14433611
$ws.Range("N132").Value = -7010
$ws.Range("H138").Value = 1891.22
$ws.Range("I138").Value = 1327.4412
$ws.Range("J138").Value = 3089.25
$ws.Range("K138").Value = 3982.3236
$ws.Range("L138").Value = 9267.75
$ws.Range("M138").Value = 1157.6764
$ws.Range("N138").Value = -19547.75
$ws.Range("H139").Value = 49663
$ws.Range("J139").Value = 49494.5
$ws.Range("L139").Value = 49494.5
$ws.Range("N139").Value = -59774.5
$ws.Range("H141").Value = 1660.4286
$ws.Range("I141").Value = 1645.6154
$ws.Range("J141").Value = 1853
$ws.Range("K141").Value = 4936.8462
$ws.Range("L141").Value = 5559
$ws.Range("M141").Value = 243.1538
$ws.Range("N141").Value = -15919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 22044.572
$ws.Range("I28").Value = 18942.4
$ws.Range("K28").Value = 18942.4
$ws.Range("M28").Value = -18750.4
$ws.Range("H61").Value = 1738.6383
$ws.Range("I61").Value = 1083.96
$ws.Range("J61").Value = 2482.5908
$ws.Range("K61").Value = 1083.96
$ws.Range("L61").Value = 2482.5908
$ws.Range("M61").Value = -871.96
$ws.Range("N61").Value = -2906.5908
$ws.Range("H99").Value = 22044.572
$ws.Range("I99").Value = 18942.4
$ws.Range("K99").Value = 18942.4
$ws.Range("M99").Value = -15947.4
$ws.Range("H124").Value = 24904.834
$ws.Range("J124").Value = 24904.834
$ws.Range("L124").Value = 24904.834
$ws.Range("N124").Value = -34724.834
$ws.Range("H125").Value = 48000
$ws.Range("J125").Value = 48000
$ws.Range("L125").Value = 48000
$ws.Range("N125").Value = -57840
$ws.Range("H132").Value = 2206.8933
$ws.Range("I132").Value = 2239.9558
$ws.Range("J132").Value = 1885.7142
$ws.Range("K132").Value = 6719.867400000001
$ws.Range("L132").Value = 5657.142599999999
$ws.Range("M132").Value = -4189.867400000001
$ws.Range("N132").Value = -10717.1426
$ws.Range("H133").Value = 38287.25
$ws.Range("J133").Value = 38287.25
$ws.Range("L133").Value = 38287.25
$ws.Range("N133").Value = -43347.25
$ws.Range("H136").Value = 1738.6383
$ws.Range("I136").Value = 1083.96
$ws.Range("J136").Value = 2482.5908
$ws.Range("K136").Value = 3251.88
$ws.Range("L136").Value = 7447.7724
$ws.Range("M136").Value = -701.8800000000001
$ws.Range("N136").Value = -12547.7724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2109.8948
$ws.Range("I99").Value = 1734.8334
$ws.Range("J99").Value = 2283
$ws.Range("K99").Value = 1734.8334
$ws.Range("L99").Value = 2283
$ws.Range("M99").Value = -236.8334
$ws.Range("N99").Value = -5279
$ws.Range("H134").Value = 2416.0488
$ws.Range("I134").Value = 2159.4211
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 6478.263300000001
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -3943.263300000001
$ws.Range("N134").Value = -22069.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 20115.666
$ws.Range("I99").Value = 4795
$ws.Range("J99").Value = 50757
$ws.Range("K99").Value = 4795
$ws.Range("L99").Value = 50757
$ws.Range("M99").Value = -3297
$ws.Range("N99").Value = -53753
$ws.Range("H124").Value = 22663
$ws.Range("J124").Value = 22663
$ws.Range("L124").Value = 22663
$ws.Range("N124").Value = -27573
$ws.Range("H126").Value = 20115.666
$ws.Range("I126").Value = 4795
$ws.Range("J126").Value = 50757
$ws.Range("K126").Value = 14385
$ws.Range("L126").Value = 152271
$ws.Range("M126").Value = -11915
$ws.Range("N126").Value = -157211
$ws.Range("H132").Value = 37502820
$ws.Range("I132").Value = 34485410
$ws.Range("J132").Value = 45457810
$ws.Range("K132").Value = 103456230
$ws.Range("L132").Value = 136373430
$ws.Range("M132").Value = -103453700
$ws.Range("N132").Value = -136378490
$ws.Range("H134").Value = 1153.6364
$ws.Range("I134").Value = 1086.5
$ws.Range("K134").Value = 3259.5
$ws.Range("M134").Value = -724.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 347950.78
$ws.Range("I107").Value = 438.25
$ws.Range("J107").Value = 695463.3
$ws.Range("K107").Value = 1314.75
$ws.Range("L107").Value = 2086389.9
$ws.Range("M107").Value = 605.25
$ws.Range("N107").Value = -2090229.9
$ws.Range("H123").Value = 3838.375
$ws.Range("I123").Value = 1986.6666
$ws.Range("K123").Value = 5959.9998
$ws.Range("M123").Value = -3509.9998
$ws.Range("H131").Value = 8380.608
$ws.Range("J131").Value = 8655.437
$ws.Range("L131").Value = 25966.311
$ws.Range("N131").Value = -36046.311

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2535.6829
$ws.Range("I132").Value = 2544.5676
$ws.Range("J132").Value = 2453.5
$ws.Range("K132").Value = 7633.702799999999
$ws.Range("L132").Value = 7360.5
$ws.Range("M132").Value = -5103.702799999999
$ws.Range("N132").Value = -12420.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 2007.7368
$ws.Range("I122").Value = 2142.9167
$ws.Range("J122").Value = 1776
$ws.Range("K122").Value = 6428.750100000001
$ws.Range("L122").Value = 5328
$ws.Range("M122").Value = -3978.750100000001
$ws.Range("N122").Value = -10228
$ws.Range("H141").Value = 74272
$ws.Range("J141").Value = 74272
$ws.Range("L141").Value = 74272
$ws.Range("N141").Value = -84632

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1973.6428
$ws.Range("I122").Value = 1601.0625
$ws.Range("J122").Value = 2470.4167
$ws.Range("K122").Value = 4803.1875
$ws.Range("L122").Value = 7411.250100000001
$ws.Range("M122").Value = -2353.1875
$ws.Range("N122").Value = -12311.2501
$ws.Range("H137").Value = 40387.855
$ws.Range("J137").Value = 40387.855
$ws.Range("L137").Value = 40387.855
$ws.Range("N137").Value = -50587.855
$ws.Range("H141").Value = 74266
$ws.Range("J141").Value = 74266
$ws.Range("L141").Value = 74266
$ws.Range("N141").Value = -84626
